$wb = $excel.ActiveWorkbook

# Sheet ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2139.077
$ws.Range("J80").Value = 2255.818
$ws.Range("L80").Value = 6767.454000000001
$ws.Range("N80").Value = -8763.454000000002

# Sheet ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 2139.077
$ws.Range("J83").Value = 2255.818
$ws.Range("L83").Value = 20302.362
$ws.Range("N83").Value = -30286.362

# Sheet ALC row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 458.63635
$ws.Range("I96").Value = 335.33334
$ws.Range("J96").Value = 504.875
$ws.Range("K96").Value = 1006.00002
$ws.Range("L96").Value = 1514.625
$ws.Range("M96").Value = 366.9999799999999
$ws.Range("N96").Value = -4260.625

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22856.666
$ws.Range("I32").Value = 3694.5957
$ws.Range("J32").Value = 151516.28
$ws.Range("K32").Value = 3694.5957
$ws.Range("L32").Value = 151516.28
$ws.Range("M32").Value = -3407.5957
$ws.Range("N32").Value = -152090.28

# Sheet ARM row 59
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

# Sheet ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7167.7144
$ws.Range("I74").Value = 1107.6428
$ws.Range("J74").Value = 19287.857
$ws.Range("K74").Value = 1107.6428
$ws.Range("L74").Value = 19287.857
$ws.Range("M74").Value = -233.6428000000001
$ws.Range("N74").Value = -21035.857

# Sheet ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 7167.7144
$ws.Range("I77").Value = 1107.6428
$ws.Range("J77").Value = 19287.857
$ws.Range("K77").Value = 5538.214
$ws.Range("L77").Value = 96439.285
$ws.Range("M77").Value = -1170.214
$ws.Range("N77").Value = -105175.285

# Sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2867.611
$ws.Range("I132").Value = 2287.7693
$ws.Range("J132").Value = 4375.2
$ws.Range("K132").Value = 6863.3079
$ws.Range("L132").Value = 13125.6
$ws.Range("M132").Value = -4333.3079
$ws.Range("N132").Value = -18185.6

# Sheet BSM row 59
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 53200
$ws.Range("J59").Value = 53200
$ws.Range("L59").Value = 53200
$ws.Range("N59").Value = -54894

# Sheet BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2051.9583
$ws.Range("I134").Value = 1349.45
$ws.Range("K134").Value = 4048.35
$ws.Range("M134").Value = -1513.35

# Sheet CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3484.8333
$ws.Range("I31").Value = 1230.2258
$ws.Range("J31").Value = 6523.6523
$ws.Range("K31").Value = 1230.2258
$ws.Range("L31").Value = 6523.6523
$ws.Range("M31").Value = -935.2257999999999
$ws.Range("N31").Value = -7113.6523

# Sheet CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3484.8333
$ws.Range("I34").Value = 1230.2258
$ws.Range("J34").Value = 6523.6523
$ws.Range("K34").Value = 1230.2258
$ws.Range("L34").Value = 6523.6523
$ws.Range("M34").Value = -1028.2258
$ws.Range("N34").Value = -6927.6523

# Sheet CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2884.96
$ws.Range("I58").Value = 1401.9231
$ws.Range("J58").Value = 4491.5835
$ws.Range("K58").Value = 1401.9231
$ws.Range("L58").Value = 4491.5835
$ws.Range("M58").Value = -1198.9231
$ws.Range("N58").Value = -4897.5835

# Sheet CRP row 68
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 22000
$ws.Range("J68").Value = 22000
$ws.Range("L68").Value = 22000
$ws.Range("N68").Value = -23498

# Sheet CRP row 71
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 22000
$ws.Range("J71").Value = 22000
$ws.Range("L71").Value = 66000
$ws.Range("N71").Value = -73488

# Sheet CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6954685.5
$ws.Range("I99").Value = 8941341
$ws.Range("J99").Value = 1390
$ws.Range("K99").Value = 8941341
$ws.Range("L99").Value = 1390
$ws.Range("M99").Value = -8939843
$ws.Range("N99").Value = -4386

# Sheet CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 361.45834
$ws.Range("I107").Value = 209.52942
$ws.Range("J107").Value = 730.4286
$ws.Range("K107").Value = 209.52942
$ws.Range("L107").Value = 730.4286
$ws.Range("M107").Value = 1710.47058
$ws.Range("N107").Value = -4570.4286

# Sheet CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 6954685.5
$ws.Range("I126").Value = 8941341
$ws.Range("J126").Value = 1390
$ws.Range("K126").Value = 26824023
$ws.Range("L126").Value = 4170
$ws.Range("M126").Value = -26821553
$ws.Range("N126").Value = -9110

# Sheet CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2569.1282
$ws.Range("I132").Value = 2334.1292
$ws.Range("J132").Value = 3479.75
$ws.Range("K132").Value = 7002.3876
$ws.Range("L132").Value = 10439.25
$ws.Range("M132").Value = -4472.3876
$ws.Range("N132").Value = -15499.25

# Sheet CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2139.275
$ws.Range("I134").Value = 908.0357
$ws.Range("J134").Value = 5012.1665
$ws.Range("K134").Value = 2724.1071
$ws.Range("L134").Value = 15036.4995
$ws.Range("M134").Value = -189.1071000000002
$ws.Range("N134").Value = -20106.4995

# Sheet CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2884.96
$ws.Range("I136").Value = 1401.9231
$ws.Range("J136").Value = 4491.5835
$ws.Range("K136").Value = 4205.7693
$ws.Range("L136").Value = 13474.7505
$ws.Range("M136").Value = -1655.7693
$ws.Range("N136").Value = -18574.7505

# Sheet CUL row 59
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

# Sheet CUL row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2998
$ws.Range("I75").Value = 1995
$ws.Range("J75").Value = 3666.6667
$ws.Range("K75").Value = 5985
$ws.Range("L75").Value = 11000.0001
$ws.Range("M75").Value = -4987
$ws.Range("N75").Value = -12996.0001

# Sheet CUL row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 2998
$ws.Range("I78").Value = 1995
$ws.Range("J78").Value = 3666.6667
$ws.Range("K78").Value = 17955
$ws.Range("L78").Value = 33000.0003
$ws.Range("M78").Value = -12963
$ws.Range("N78").Value = -42984.0003

# Sheet CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5651034
$ws.Range("J131").Value = 6174256
$ws.Range("L131").Value = 18522768
$ws.Range("N131").Value = -18532848

# Sheet CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 5021.515
$ws.Range("I140").Value = 6795
$ws.Range("J140").Value = 2893.3333
$ws.Range("K140").Value = 20385
$ws.Range("L140").Value = 8679.999899999999
$ws.Range("M140").Value = -15205
$ws.Range("N140").Value = -19039.9999

# Sheet GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3047.652
$ws.Range("I80").Value = 2994.8
$ws.Range("J80").Value = 3400
$ws.Range("K80").Value = 2994.8
$ws.Range("L80").Value = 3400
$ws.Range("M80").Value = -1996.8
$ws.Range("N80").Value = -5396

# Sheet GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3047.652
$ws.Range("I83").Value = 2994.8
$ws.Range("J83").Value = 3400
$ws.Range("K83").Value = 14974
$ws.Range("L83").Value = 17000
$ws.Range("M83").Value = -9982
$ws.Range("N83").Value = -26984

# Sheet GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 840
$ws.Range("I97").Value = 840
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 840
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -344
$ws.Range("N97").ClearContents()

# Sheet GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 928762.2
$ws.Range("I122").Value = 2779029.5
$ws.Range("K122").Value = 8337088.5
$ws.Range("M122").Value = -8334638.5

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2532.2424
$ws.Range("I132").Value = 2534.08
$ws.Range("J132").Value = 2526.5
$ws.Range("K132").Value = 7602.24
$ws.Range("L132").Value = 7579.5
$ws.Range("M132").Value = -5072.24
$ws.Range("N132").Value = -12639.5

# Sheet LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12033.556
$ws.Range("I22").Value = 800.3333
$ws.Range("J22").Value = 34500
$ws.Range("K22").Value = 800.3333
$ws.Range("L22").Value = 34500
$ws.Range("M22").Value = -505.3333
$ws.Range("N22").Value = -35090

# Sheet LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 12033.556
$ws.Range("I27").Value = 800.3333
$ws.Range("J27").Value = 34500
$ws.Range("K27").Value = 800.3333
$ws.Range("L27").Value = 34500
$ws.Range("M27").Value = -693.3333
$ws.Range("N27").Value = -34714

# Sheet LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2622
$ws.Range("I68").Value = 2500
$ws.Range("J68").Value = 2774.5
$ws.Range("K68").Value = 2500
$ws.Range("L68").Value = 2774.5
$ws.Range("M68").Value = -1751
$ws.Range("N68").Value = -4272.5

# Sheet LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2622
$ws.Range("I71").Value = 2500
$ws.Range("J71").Value = 2774.5
$ws.Range("K71").Value = 12500
$ws.Range("L71").Value = 13872.5
$ws.Range("M71").Value = -8756
$ws.Range("N71").Value = -21360.5

# Sheet LTW row 116
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 30000
$ws.Range("J116").Value = 30000
$ws.Range("L116").Value = 30000
$ws.Range("N116").Value = -39178

# Sheet WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 266.66666
$ws.Range("I107").Value = 266.66666
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 799.9999799999999
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1120.00002
$ws.Range("N107").ClearContents()

# Sheet WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 527.8889
$ws.Range("I113").Value = 403.66666
$ws.Range("J113").Value = 776.3333
$ws.Range("K113").Value = 1210.99998
$ws.Range("L113").Value = 2328.9999
$ws.Range("M113").Value = 959.0000199999999
$ws.Range("N113").Value = -6668.9999

# Sheet WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 101147.3
$ws.Range("I122").Value = 125871.625
$ws.Range("K122").Value = 377614.875
$ws.Range("M122").Value = -375164.875

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3226.675
$ws.Range("I132").Value = 3359.9678
$ws.Range("J132").Value = 2767.5557
$ws.Range("K132").Value = 10079.9034
$ws.Range("L132").Value = 8302.667099999999
$ws.Range("M132").Value = -7549.903399999999
$ws.Range("N132").Value = -13362.6671

# Sheet WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3167.9167
$ws.Range("I136").Value = 1541.3889
$ws.Range("J136").Value = 4794.4443
$ws.Range("K136").Value = 4624.1667
$ws.Range("L136").Value = 14383.3329
$ws.Range("M136").Value = -2074.1667
$ws.Range("N136").Value = -19483.3329
